# Apply the StructureDefinition-vaccine-type.xlsx metadata update:
#  1. Update the URL (row 2 / column B) from the "pythia" IG to "cicada".
#  2. Update the Date value (row 8 / column B) to the new generation timestamp.
#  3. Insert a new "Jurisdiction" property row (with an empty value) right
#     after "Contact" and before "Description" on the Metadata sheet, shifting
#     every following row down by one.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")
$ws2 = $wb.Worksheets.Item("Elements")

# -- 1. URL --------------------------------------------------------------
# The canonical IG URL moved from the "pythia" to the "cicada" IG. The same
# string is also used as the Fixed Value of the Extension.url element on the
# Elements sheet (column R, row 5), so update both occurrences to keep them
# sharing one string, exactly as before the edit.
$newUrl = "http://fhirfli.dev/fhir/ig/cicada/StructureDefinition/vaccine-type"
$ws1.Cells.Item(2, 2).Value = $newUrl
$ws2.Cells.Item(5, 18).Value = $newUrl

# -- 2. Date ---------------------------------------------------------------
$ws1.Cells.Item(8, 2).Value = "2026-02-11T14:37:07-05:00"

# -- 3. Insert the Jurisdiction row ----------------------------------------
# Insert a blank row at row 11 (pushing "Description" and everything after
# it down by one row), then copy the formatting from the row that is now
# directly below (the old "Contact" + 2 -> now "Description" row, which
# still carries the same border/alignment style used throughout the table)
# so the new row matches the existing look instead of picking up a blank
# default style.
$ws1.Rows.Item(11).Insert()
$ws1.Range("A12:B12").Copy()
$ws1.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws1.Cells.Item(11, 1).Value = "Jurisdiction"
$ws1.Cells.Item(11, 2).Value = ""
